# PlayerPerformance_4590.xlsx — add "ODI Batting Extra" / "ODI Bowling Extra"
# sheets, and clear the stray empty B3/B7 cells on "ODI Batting".

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$value) {
    # Force the cell to be stored as text (matches the workbook's existing
    # inline-string cells) without leaving a lingering custom number format
    # / style index behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

function Format-HeaderRow($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
    $range.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# 1. "ODI Batting": B3 / B7 used to hold an empty inline string; the
#    cells should simply not exist any more.
# ---------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("B3").ClearContents()
$battingSheet.Range("B7").ClearContents()

# ---------------------------------------------------------------------
# 2. New sheet "ODI Batting Extra" (after "ODI Bowling")
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$battingExtraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 0; $col -lt $battingExtraHeaders.Length; $col++) {
    $battingExtra.Cells.Item(1, $col + 1).Value = $battingExtraHeaders[$col]
}
Format-HeaderRow $battingExtra.Range("A1:F1")

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# BATTING_POSITION ($null) is numeric where present; everything else is text.
$battingExtraRows = @(
    @("3936", $null, $null, $null, $null,   "NO"),
    @("3938", 10,    "0",   "0",   "0.52%", "NO"),
    @("4377", 11,    "0",   "0",   "0.70%", "NO"),
    @("4378", 11,    "2",   "0",   "3.76%", "NO"),
    @("4444", $null, $null, $null, $null,   $null),
    @("4446", $null, $null, $null, $null,   $null),
    @("4448", $null, $null, $null, $null,   $null)
)

for ($i = 0; $i -lt $battingExtraRows.Length; $i++) {
    $row = 2 + $i
    $data = $battingExtraRows[$i]

    Set-TextValue $battingExtra.Cells.Item($row, 1) $data[0]

    if ($null -ne $data[1]) {
        $battingExtra.Cells.Item($row, 2).Value = $data[1]
    }
    if ($null -ne $data[2]) {
        Set-TextValue $battingExtra.Cells.Item($row, 3) $data[2]
    }
    if ($null -ne $data[3]) {
        Set-TextValue $battingExtra.Cells.Item($row, 4) $data[3]
    }
    if ($null -ne $data[4]) {
        Set-TextValue $battingExtra.Cells.Item($row, 5) $data[4]
    }
    if ($null -ne $data[5]) {
        Set-TextValue $battingExtra.Cells.Item($row, 6) $data[5]
    }
}

# ---------------------------------------------------------------------
# 3. New sheet "ODI Bowling Extra" (after "ODI Batting Extra")
# ---------------------------------------------------------------------
$bowlingExtra = $wb.Worksheets.Add($null, $battingExtra)
$bowlingExtra.Name = "ODI Bowling Extra"

$bowlingExtraHeaders = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 0; $col -lt $bowlingExtraHeaders.Length; $col++) {
    $bowlingExtra.Cells.Item(1, $col + 1).Value = $bowlingExtraHeaders[$col]
}
Format-HeaderRow $bowlingExtra.Range("A1:C1")

$bowlingExtraRows = @(
    @("3936", "0",   "10.00%"),
    @("3938", $null, $null),
    @("4377", "1",   "10.00%"),
    @("4378", "0",   "30.00%"),
    @("4444", "0",   "30.00%"),
    @("4446", "0",   "40.00%"),
    @("4448", "0",   "10.00%")
)

for ($i = 0; $i -lt $bowlingExtraRows.Length; $i++) {
    $row = 2 + $i
    $data = $bowlingExtraRows[$i]

    Set-TextValue $bowlingExtra.Cells.Item($row, 1) $data[0]

    if ($null -ne $data[1]) {
        Set-TextValue $bowlingExtra.Cells.Item($row, 2) $data[1]
    }
    if ($null -ne $data[2]) {
        Set-TextValue $bowlingExtra.Cells.Item($row, 3) $data[2]
    }
}
